$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 479.5
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()

$ws.Range("H62").Value = 2999
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()

$ws.Range("H65").Value = 2999
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()

$ws.Range("H80").Value = 1717.875
$ws.Range("I80").Value = 1068.2858
$ws.Range("J80").Value = 1985.3529
$ws.Range("K80").Value = 3204.8574
$ws.Range("L80").Value = 5956.0587
$ws.Range("M80").Value = -2206.8574
$ws.Range("N80").Value = -7952.0587

$ws.Range("H83").Value = 1717.875
$ws.Range("I83").Value = 1068.2858
$ws.Range("J83").Value = 1985.3529
$ws.Range("K83").Value = 9614.572200000001
$ws.Range("L83").Value = 17868.1761
$ws.Range("M83").Value = -4622.572200000001
$ws.Range("N83").Value = -27852.1761

$ws.Range("H103").Value = 1574.8334
$ws.Range("I103").Value = 1162.25
$ws.Range("J103").Value = 2400
$ws.Range("K103").Value = 3486.75
$ws.Range("L103").Value = 7200
$ws.Range("M103").Value = -2900.75
$ws.Range("N103").Value = -8372

$ws.Range("H106").Value = 2531.3333
$ws.Range("I106").Value = 2997
$ws.Range("J106").Value = 1600
$ws.Range("K106").Value = 2997
$ws.Range("L106").Value = 1600
$ws.Range("M106").Value = -2366
$ws.Range("N106").Value = -2862

$ws.Range("H113").Value = 5323.5713
$ws.Range("I113").Value = 5942.5
$ws.Range("K113").Value = 5942.5
$ws.Range("M113").Value = -2688.5

$ws.Range("H138").Value = 4997.5835
$ws.Range("J138").Value = 8439.842000000001
$ws.Range("L138").Value = 25319.526
$ws.Range("N138").Value = -35599.526

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5599.3613
$ws.Range("I32").Value = 5187.914
$ws.Range("K32").Value = 5187.914
$ws.Range("M32").Value = -4900.914

$ws.Range("H45").Value = 3272.9583
$ws.Range("I45").Value = 2930.9
$ws.Range("J45").Value = 3517.2856
$ws.Range("K45").Value = 2930.9
$ws.Range("L45").Value = 3517.2856
$ws.Range("M45").Value = -2553.9
$ws.Range("N45").Value = -4271.2856

$ws.Range("H74").Value = 27229.355
$ws.Range("I74").Value = 27229.355
$ws.Range("K74").Value = 27229.355
$ws.Range("M74").Value = -26355.355

$ws.Range("H77").Value = 27229.355
$ws.Range("I77").Value = 27229.355
$ws.Range("K77").Value = 136146.775
$ws.Range("M77").Value = -131778.775

$ws.Range("H97").Value = 1040.2778
$ws.Range("I97").Value = 898.86664
$ws.Range("K97").Value = 898.86664
$ws.Range("M97").Value = -402.86664

$ws.Range("H130").Value = 60547.832
$ws.Range("J130").Value = 60547.832
$ws.Range("L130").Value = 60547.832
$ws.Range("N130").Value = -70587.83199999999

$ws.Range("H131").Value = 53204.332
$ws.Range("J131").Value = 53204.332
$ws.Range("L131").Value = 53204.332
$ws.Range("N131").Value = -63284.332

$ws.Range("H132").Value = 107756.07
$ws.Range("I132").Value = 8807.083000000001
$ws.Range("K132").Value = 26421.249
$ws.Range("M132").Value = -23891.249

$ws.Range("H139").Value = 95282.86
$ws.Range("J139").Value = 95282.86
$ws.Range("L139").Value = 95282.86
$ws.Range("N139").Value = -105562.86

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H16").Value = 1000
$ws.Range("I16").Value = 1000
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 1000
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -830
$ws.Range("N16").ClearContents()

$ws.Range("H105").Value = 5618.528
$ws.Range("I105").Value = 4025
$ws.Range("K105").Value = 4025
$ws.Range("M105").Value = -2278

$ws.Range("H134").Value = 2085.111
$ws.Range("I134").Value = 2095.75
$ws.Range("K134").Value = 6287.25
$ws.Range("M134").Value = -3752.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 5630.7
$ws.Range("I22").Value = 6557.125
$ws.Range("J22").Value = 1925
$ws.Range("K22").Value = 6557.125
$ws.Range("L22").Value = 1925
$ws.Range("M22").Value = -6207.125
$ws.Range("N22").Value = -2625

$ws.Range("H31").Value = 8336680
$ws.Range("I31").Value = 2577.6191
$ws.Range("K31").Value = 2577.6191
$ws.Range("M31").Value = -2282.6191

$ws.Range("H34").Value = 8336680
$ws.Range("I34").Value = 2577.6191
$ws.Range("K34").Value = 2577.6191
$ws.Range("M34").Value = -2375.6191

$ws.Range("H131").Value = 34974.75
$ws.Range("J131").Value = 34974.75
$ws.Range("L131").Value = 34974.75
$ws.Range("N131").Value = -45054.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 998.3333
$ws.Range("I68").Value = 995
$ws.Range("K68").Value = 2985
$ws.Range("M68").Value = -2174

$ws.Range("H71").Value = 998.3333
$ws.Range("I71").Value = 995
$ws.Range("K71").Value = 8955
$ws.Range("M71").Value = -4899

$ws.Range("H80").Value = 3951.5
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 3951.5
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 11854.5
$ws.Range("N80").Value = -13726.5
$ws.Range("M80").ClearContents()

$ws.Range("H83").Value = 3951.5
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 3951.5
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 35563.5
$ws.Range("N83").Value = -44923.5
$ws.Range("M83").ClearContents()

$ws.Range("H122").Value = 1192.1177
$ws.Range("J122").Value = 1388.5
$ws.Range("L122").Value = 12496.5
$ws.Range("N122").Value = -17396.5

$ws.Range("H132").Value = 5077.3
$ws.Range("J132").Value = 8889.799999999999
$ws.Range("L132").Value = 80008.2
$ws.Range("N132").Value = -85068.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1040.826
$ws.Range("J16").Value = 1000
$ws.Range("L16").Value = 1000
$ws.Range("N16").Value = -1340

$ws.Range("H68").Value = 9999.666999999999
$ws.Range("I68").Value = 9999.5
$ws.Range("K68").Value = 9999.5
$ws.Range("M68").Value = -9250.5

$ws.Range("H71").Value = 9999.666999999999
$ws.Range("I71").Value = 9999.5
$ws.Range("K71").Value = 49997.5
$ws.Range("M71").Value = -46253.5

$ws.Range("H100").Value = 3670.0356
$ws.Range("I100").Value = 3367.6667
$ws.Range("K100").Value = 3367.6667
$ws.Range("M100").Value = -2826.6667

$ws.Range("H132").Value = 2538.611
$ws.Range("J132").Value = 2869
$ws.Range("L132").Value = 8607
$ws.Range("N132").Value = -13667

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1279.9584
$ws.Range("I107").Value = 972
$ws.Range("J107").Value = 2819.75
$ws.Range("K107").Value = 2916
$ws.Range("L107").Value = 8459.25
$ws.Range("M107").Value = -996
$ws.Range("N107").Value = -12299.25

$ws.Range("H122").Value = 12504736
$ws.Range("I122").Value = 3856.6667
$ws.Range("K122").Value = 11570.0001
$ws.Range("M122").Value = -9120.000100000001

$ws.Range("H126").Value = 18523084
$ws.Range("I126").Value = 20837844
$ws.Range("J126").Value = 5000
$ws.Range("K126").Value = 62513532
$ws.Range("L126").Value = 15000
$ws.Range("M126").Value = -62511062
$ws.Range("N126").Value = -19940

$ws.Range("H132").Value = 2229.5356
$ws.Range("I132").Value = 2074.0625
$ws.Range("J132").Value = 3162.375
$ws.Range("K132").Value = 6222.1875
$ws.Range("L132").Value = 9487.125
$ws.Range("M132").Value = -3692.1875
$ws.Range("N132").Value = -14547.125

$ws.Range("H136").Value = 3981.5264
$ws.Range("I136").Value = 2790.625
$ws.Range("J136").Value = 10333
$ws.Range("K136").Value = 8371.875
$ws.Range("L136").Value = 30999
$ws.Range("M136").Value = -5821.875
$ws.Range("N136").Value = -36099
